$wb = $excel.ActiveWorkbook
$existing = $wb.Worksheets.Item(1)
$existing.Copy($null, $existing)
$copy = $wb.Worksheets.Item(2)
$existing.Name = "Limón automatizado"
$copy.Name = "Limón manual"

# Edits specific to "Limón automatizado" (sheet1 / $existing)
$existing.Range("A27").Value = "selección por tamaño"
$existing.Range("B27").Value = 1
$existing.Range("B33").Value = 3
$existing.Range("B36").Value = 3

# Edits specific to "Limón manual" (sheet2 / $copy)
$copy.Range("B27").Value = 7

$existing.Activate()
